$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The original sheet had 6 data rows (rows 2-7), all "FAPs -> <target>" pairs.
# This edit adds 6 new "ECs -> <target>" rows ahead of them, pushing the
# original FAPs rows down to rows 8-13, and refreshes every numeric
# statistic now that the extra source cluster is part of the model.
#
# Grow the sheet by inserting the 6 extra rows at the bottom of the
# existing data block (after row 7) rather than right under the bold/
# bordered header row - that keeps the new rows from inheriting the
# header's style. We then overwrite every cell's content below, so it
# does not matter that the "new" rows are appended rather than inserted
# at row 2 - the final row-by-row values are written explicitly anyway.
$ws.Range("A8:A13").EntireRow.Insert()
$ws.Range("A8:T13").ClearFormats()

$allRows = @(
    @("ECs","Wnt2","Fzd1","ECs",1,0.3333333333333333,0.01070233333333333,0.032107,0.006017198313602724,0.006017198313602724,3,1,0.5018866666666667,1.50566,0.02325346250969519,0.02325346250969519,0.005371358402222222,0.04834222562000001,0.0001399206953987621,0.0001399206953987621),
    @("ECs","Wnt2","Fzd1","FAPs",1,0.3333333333333333,0.01070233333333333,0.032107,0.006017198313602724,0.006017198313602724,3,1,17.00656966666667,51.019709,0.7879500620904175,0.7879500620904175,0.1820099774292222,1.638089796863,0.004741251784813622,0.004741251784813622),
    @("ECs","Wnt2","Fzd1","M1",1,0.3333333333333333,0.01070233333333333,0.032107,0.006017198313602724,0.006017198313602724,3,1,0.4244703333333333,1.273411,0.01966660132296365,0.01966660132296365,0.004542822997444444,0.04088540697699999,0.000118337840314834,0.000118337840314834),
    @("ECs","Wnt2","Fzd1","M2",1,0.3333333333333333,0.01070233333333333,0.032107,0.006017198313602724,0.006017198313602724,3,1,0.3996003333333333,1.198801,0.01851432203159086,0.01851432203159086,0.004276655967444444,0.038489903707,0.0001114043473059863,0.0001114043473059863),
    @("ECs","Wnt2","Fzd1","Neutro",1,0.3333333333333333,0.01070233333333333,0.032107,0.006017198313602724,0.006017198313602724,3,1,0.2431653333333333,0.7294959999999999,0.01126636019218987,0.01126636019218987,0.002602436452444444,0.02342192807199999,0.00006779192354888574,0.00006779192354888574),
    @("ECs","Wnt2","Fzd1","sCs",1,0.3333333333333333,0.01070233333333333,0.032107,0.006017198313602724,0.006017198313602724,3,1,3.007616666666667,9.02285,0.139349191853143,0.1393491918531429,0.03218851610555556,0.28969664495,0.0008384917222206342,0.0008384917222206341),
    @("FAPs","Wnt2","Fzd1","ECs",3,1,1.767921666666667,5.303765,0.9939828016863973,0.9939828016863973,3,1,0.5018866666666667,1.50566,0.02325346250969519,0.02325346250969519,0.8872963122111113,7.985666809900001,0.02311354181429643,0.02311354181429643),
    @("FAPs","Wnt2","Fzd1","FAPs",3,1,1.767921666666667,5.303765,0.9939828016863973,0.9939828016863973,3,1,17.00656966666667,51.019709,0.7879500620904175,0.7879500620904175,30.06628298937612,270.596546904385,0.783208810305604,0.783208810305604),
    @("FAPs","Wnt2","Fzd1","M1",3,1,1.767921666666667,5.303765,0.9939828016863973,0.9939828016863973,3,1,0.4244703333333333,1.273411,0.01966660132296365,0.01966660132296365,0.7504302991572223,6.753872692414999,0.01954826348264882,0.01954826348264882),
    @("FAPs","Wnt2","Fzd1","M2",3,1,1.767921666666667,5.303765,0.9939828016863973,0.9939828016863973,3,1,0.3996003333333333,1.198801,0.01851432203159086,0.01851432203159086,0.7064620873072223,6.358158785765,0.01840291768428488,0.01840291768428488),
    @("FAPs","Wnt2","Fzd1","Neutro",3,1,1.767921666666667,5.303765,0.9939828016863973,0.9939828016863973,3,1,0.2431653333333333,0.7294959999999999,0.01126636019218987,0.01126636019218987,0.4298972613822222,3.86907535244,0.01119856826864098,0.01119856826864098),
    @("FAPs","Wnt2","Fzd1","sCs",3,1,1.767921666666667,5.303765,0.9939828016863973,0.9939828016863973,3,1,3.007616666666667,9.02285,0.139349191853143,0.1393491918531429,5.317230670027778,47.85507603025,0.1385107001309223,0.1385107001309223)
)

for ($i = 0; $i -lt $allRows.Length; $i++) {
    $r = 2 + $i
    $data = $allRows[$i]
    for ($c = 1; $c -le 20; $c++) {
        $ws.Cells.Item($r, $c).Value = $data[$c - 1]
    }
}
